# This script applies a cyclic rotation of the per-observation data
# (date, quality, volume, prices, commercialization unit, etc.) across
# rows 2-14 of the sheet, as described by the commit "Fruta / hortaliza,
# semanal": each row's weekly-observation fields are replaced by those
# that used to belong to another row, per a fixed permutation derived
# from the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values move together as one "observation" record.
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

# Row permutation: destination row -> source row (values read from the
# source row in the original workbook are written into the destination
# row).
$mapping = @{
    2  = 6
    3  = 7
    4  = 11
    5  = 12
    6  = 4
    7  = 5
    8  = 2
    9  = 3
    10 = 13
    11 = 14
    12 = 8
    13 = 9
    14 = 10
}

# First, snapshot the current ("before") values for every relevant cell,
# so that writing results doesn't clobber values we still need to read.
$snapshot = @{}
foreach ($row in 2..14) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $snapshot[$addr] = $ws.Range($addr).Value2
    }
}

# Now write the rotated values back into each destination row.
foreach ($destRow in 2..14) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $srcAddr = "$col$srcRow"
        $destAddr = "$col$destRow"
        $ws.Range($destAddr).Value = $snapshot[$srcAddr]
    }
}
